# Append the newest daily cumulative Covid death figures (rows 399-410,
# dates 2021-11-18 .. 2021-11-29) to the "Deaths" sheet, matching the
# "Add files via upload" refresh of the source data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows: Date (serial), DeathCovid, DeathWithCovid, Total
$newRows = @(
    @(44518, 13781, 2629, 16410),
    @(44519, 13818, 2638, 16456),
    @(44520, 13861, 2638, 16499),
    @(44521, 13919, 2638, 16557),
    @(44522, 13985, 2649, 16634),
    @(44523, 14056, 2653, 16709),
    @(44524, 14107, 2655, 16762),
    @(44525, 14177, 2662, 16839),
    @(44526, 14228, 2670, 16898),
    @(44527, 14274, 2678, 16952),
    @(44528, 14341, 2687, 17028),
    @(44529, 14418, 2696, 17114)
)

$firstRow = 399
$lastRow = $firstRow + $newRows.Count - 1

# Build a 2D array and write it in one shot, like pasting a block of data.
$block = New-Object 'object[,]' $newRows.Count,4
for ($i = 0; $i -lt $newRows.Count; $i++) {
    for ($j = 0; $j -lt 4; $j++) {
        $block[$i, $j] = $newRows[$i][$j]
    }
}

$targetRange = $ws.Range("A$firstRow`:D$lastRow")
$targetRange.Value = $block

# Move the selection down to the newly-added last row and scroll the
# window so that row is visible, matching the saved view state
# (topLeftCell="A364", selection activeCell="A410").
$lastCell = $ws.Cells.Item($lastRow, 1)
$lastCell.Select()
$excel.ActiveWindow.ScrollRow = 364
$excel.ActiveWindow.ScrollColumn = 1
